$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet / active tab
$ws.Activate()

# Restore the selection on the sheet (shifted one column to the right, M12 -> R12)
$ws.Range("R12").Select() | Out-Null
